$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New header row styling: bold font applied to the four existing headers
# (A1:D1) plus three new header blocks (F1:K1, M1:R1). Doing this first
# creates font index 2 ("bold") / cellXfs index 3 in the style tables.
# ---------------------------------------------------------------------------
$ws.Range("A1:D1").Font.Bold = $true

# ---------------------------------------------------------------------------
# Second header block: F1:K1 -> kota / penduduk / perokok / tb / probability
# of selection / probability of inclusion
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "kota"
$ws.Range("G1").Value = "penduduk"
$ws.Range("H1").Value = "perokok"
$ws.Range("I1").Value = "tb"
$ws.Range("J1").Value = "probability of selection"
$ws.Range("K1").Value = "probability of inclusion"
$ws.Range("F1:K1").Font.Bold = $true

# ---------------------------------------------------------------------------
# Third header block: M1:P1 same headers again, Q1:R1 blank (but still
# bold-styled, matching the header row format).
# ---------------------------------------------------------------------------
$ws.Range("M1").Value = "kota"
$ws.Range("N1").Value = "penduduk"
$ws.Range("O1").Value = "perokok"
$ws.Range("P1").Value = "tb"
$ws.Range("M1:R1").Font.Bold = $true

# ---------------------------------------------------------------------------
# Sample table (F2:K5): 4 of the 6 cities with computed sampling
# probabilities (probability of selection = penduduk / total penduduk;
# probability of inclusion from a PPS draw-by-draw simulation).
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = "JAKARTA BARAT"
$ws.Range("G2").Value = 2569462
$ws.Range("H2").Value = 387218.74249999999
$ws.Range("I2").Value = 4462
$ws.Range("J2").Value = 0.229485238999974
$ws.Range("K2").Value = 0.64749999999999996

$ws.Range("F3").Value = "JAKARTA PUSAT"
$ws.Range("G3").Value = 1153399
$ws.Range("H3").Value = 165788.98510000002
$ws.Range("I3").Value = 6670
$ws.Range("J3").Value = 0.103013021861125
$ws.Range("K3").Value = 0.35260000000000002

$ws.Range("F4").Value = "JAKARTA TIMUR"
$ws.Range("G4").Value = 3234003
$ws.Range("H4").Value = 474220.0932
$ws.Range("I4").Value = 4126
$ws.Range("J4").Value = 0.28883709951018299
$ws.Range("K4").Value = 0.74419999999999997

$ws.Range("F5").Value = "JAKARTA UTARA"
$ws.Range("G5").Value = 1843537
$ws.Range("H5").Value = 273636.40590000001
$ws.Range("I5").Value = 2662
$ws.Range("J5").Value = 0.16465101606884799
$ws.Range("K5").Value = 0.5131

# Match the number formatting of the original table: penduduk/perokok use
# the thousands-style numeric format, tb uses the "Normal 2" cell style
# (same as the source table's D column).
$ws.Range("G2:H5").NumberFormat = "#,##0"
$ws.Range("I2:I5").Style = "Normal 2"

# Probability columns get a plain-black font (new font/cellXf entries).
$ws.Range("J2:K5").Font.Color = 0

# ---------------------------------------------------------------------------
# Column L stays empty in every sampled row, but still carries its own
# (slightly different) black font style - mirrors the source file's
# leftover formatting-only paste.
# ---------------------------------------------------------------------------
$ws.Range("L2:L5").Font.Color = 0
$ws.Range("L2:L5").Font.Name = "Calibri "

# ---------------------------------------------------------------------------
# Final subset table (M2:P4): three of the four sampled cities, re-pasted
# without the probability columns. Q2:R4 stay empty but keep the
# plain-black font used for the probability columns.
# ---------------------------------------------------------------------------
$ws.Range("M2").Value = "JAKARTA BARAT"
$ws.Range("N2").Value = 2569462
$ws.Range("O2").Value = 387218.74249999999
$ws.Range("P2").Value = 4462

$ws.Range("M3").Value = "JAKARTA TIMUR"
$ws.Range("N3").Value = 3234003
$ws.Range("O3").Value = 474220.0932
$ws.Range("P3").Value = 4126

$ws.Range("M4").Value = "JAKARTA UTARA"
$ws.Range("N4").Value = 1843537
$ws.Range("O4").Value = 273636.40590000001
$ws.Range("P4").Value = 2662

$ws.Range("N2:O4").NumberFormat = "#,##0"
$ws.Range("P2:P4").Style = "Normal 2"

$ws.Range("Q2:R4").Font.Color = 0

# ---------------------------------------------------------------------------
# Column widths for the new blocks (best-effort; engine quantizes these).
# ---------------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 15.33
$ws.Columns("G").ColumnWidth = 10.11
$ws.Columns("J").ColumnWidth = 19.78
$ws.Columns("K").ColumnWidth = 19.66
$ws.Columns("M").ColumnWidth = 18.11
$ws.Columns("N").ColumnWidth = 11.22

# ---------------------------------------------------------------------------
# View state: scrolled so column K is first visible, L7 selected.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 11
$ws.Range("L7").Select() | Out-Null

# ---------------------------------------------------------------------------
# Page setup: portrait orientation (print-quality DPI / printer-settings
# relationship are not exposed on this COM surface).
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
